$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in table values for ANN (column C) & SVM (column E) for rows 5 and 6
$ws.Range("C5").Value = 0.74
$ws.Range("E5").Value = 0.2893

$ws.Range("C6").Value = 0.759
$ws.Range("E6").Value = 0.0793

# Copy number formatting/style from the row above (C4/E4) so the new
# percentage cells match the existing table styling
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C5:C6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("E4").Copy() | Out-Null
$ws.Range("E5:E6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# Match styling on the blank D/F cells in these rows to the rest of the table
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D5:D6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("F4").Copy() | Out-Null
$ws.Range("F5:F6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# Update the active selection to match the edit
$ws.Range("C7").Select() | Out-Null
